$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 72
$ws.Range("E2").Value = 3041.325569152832
$ws.Range("F2").Value = 20

$ws.Range("C3").Value = 68
$ws.Range("D3").Value = 4
$ws.Range("E3").Value = 2740.487575531006
$ws.Range("F3").Value = 26
$ws.Range("G3").Value = 16

$ws.Rows("4:11").Delete()
